$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.2992651321056314
$ws.Cells.Item(2, 3).Value = 0.03537740902002895
$ws.Cells.Item(2, 4).Value = 0.0302029487070179
$ws.Cells.Item(2, 6).Value = 0.911962873044402
$ws.Cells.Item(2, 7).Value = 0.7620888783238087
$ws.Cells.Item(2, 8).Value = 0.8275719869384233
$ws.Cells.Item(2, 11).Value = 0.2580512119084233
$ws.Cells.Item(2, 13).Value = 0.8842635028906898
$ws.Cells.Item(2, 14).Value = 1.698872570602802

$ws.Cells.Item(3, 2).Value = 0.2693374846342635
$ws.Cells.Item(3, 3).Value = 0.03226038484955041
$ws.Cells.Item(3, 4).Value = 0.02925759992202615
$ws.Cells.Item(3, 6).Value = 0.9029426628173383
$ws.Cells.Item(3, 7).Value = 0.7542035276183725
$ws.Cells.Item(3, 8).Value = 0.8281616380371162
$ws.Cells.Item(3, 11).Value = 0.2276848303391574
$ws.Cells.Item(3, 13).Value = 0.7874571915405966
$ws.Cells.Item(3, 14).Value = 1.713920251458752

$ws.Cells.Item(4, 2).Value = 0.2510503727820605
$ws.Cells.Item(4, 3).Value = 0.03032946657226887
$ws.Cells.Item(4, 4).Value = 0.02866764449029091
$ws.Cells.Item(4, 6).Value = 0.8979311507071799
$ws.Cells.Item(4, 7).Value = 0.7498336924256392
$ws.Cells.Item(4, 8).Value = 0.8288897596504512
$ws.Cells.Item(4, 11).Value = 0.2090696262872029
$ws.Cells.Item(4, 13).Value = 0.7285494954619338
$ws.Cells.Item(4, 14).Value = 1.723748593703235

$ws.Cells.Item(5, 2).Value = 0.2436207805673973
$ws.Cells.Item(5, 3).Value = 0.02953834078065398
$ws.Cells.Item(5, 4).Value = 0.0284248633566051
$ws.Cells.Item(5, 6).Value = 0.8960212325221875
$ws.Cells.Item(5, 7).Value = 0.7481712900979431
$ws.Cells.Item(5, 8).Value = 0.8292784951110832
$ws.Cells.Item(5, 11).Value = 0.2014915633193226
$ws.Cells.Item(5, 13).Value = 0.7046717085652574
$ws.Cells.Item(5, 14).Value = 1.727901593689253

$ws.Cells.Item(6, 2).Value = 0.2423884729192878
$ws.Cells.Item(6, 3).Value = 0.02940671839564146
$ws.Cells.Item(6, 4).Value = 0.02838440724648805
$ws.Cells.Item(6, 6).Value = 0.8957120780517798
$ws.Cells.Item(6, 7).Value = 0.7479023883238511
$ws.Cells.Item(6, 8).Value = 0.8293486012188822
$ws.Cells.Item(6, 11).Value = 0.2002337098711706
$ws.Cells.Item(6, 13).Value = 0.7007143481772857
$ws.Cells.Item(6, 14).Value = 1.728600120253859

$ws.Cells.Item(7, 2).Value = 0.2509500830099114
$ws.Cells.Item(7, 3).Value = 0.03031881438927542
$ws.Cells.Item(7, 4).Value = 0.02866437982522285
$ws.Cells.Item(7, 6).Value = 0.8979048574030912
$ws.Cells.Item(7, 7).Value = 0.7498107939264855
$ws.Cells.Item(7, 8).Value = 0.8288946297353448
$ws.Cells.Item(7, 11).Value = 0.2089673940031389
$ws.Cells.Item(7, 13).Value = 0.7282269627317106
$ws.Cells.Item(7, 14).Value = 1.723804004127182

$ws.Cells.Item(8, 2).Value = 0.2889278510871804
$ws.Cells.Item(8, 3).Value = 0.03430621386409172
$ws.Cells.Item(8, 4).Value = 0.02987897818983498
$ws.Cells.Item(8, 6).Value = 0.908743220459229
$ws.Cells.Item(8, 7).Value = 0.7592719026954029
$ws.Cells.Item(8, 8).Value = 0.8276992985167766
$ws.Cells.Item(8, 11).Value = 0.2475748356261107
$ws.Cells.Item(8, 13).Value = 0.8507705626288811
$ws.Cells.Item(8, 14).Value = 1.703938619192542

$ws.Cells.Item(9, 2).Value = 0.3640966742646299
$ws.Cells.Item(9, 3).Value = 0.04198932602071181
$ws.Cells.Item(9, 4).Value = 0.03218452301479857
$ws.Cells.Item(9, 6).Value = 0.9341891021150985
$ws.Cells.Item(9, 7).Value = 0.7815841374593049
$ws.Cells.Item(9, 8).Value = 0.8282622403939399
$ws.Cells.Item(9, 11).Value = 0.3235128745002953
$ws.Cells.Item(9, 13).Value = 1.095590066422005
$ws.Cells.Item(9, 14).Value = 1.669666600448416

$ws.Cells.Item(10, 2).Value = 0.4197411787082501
$ws.Cells.Item(10, 3).Value = 0.04755043368105305
$ws.Cells.Item(10, 4).Value = 0.03383090674083888
$ws.Cells.Item(10, 6).Value = 0.955458280623958
$ws.Cells.Item(10, 7).Value = 0.8002929908857794
$ws.Cells.Item(10, 8).Value = 0.8304525804958587
$ws.Cells.Item(10, 11).Value = 0.3794388716311801
$ws.Cells.Item(10, 13).Value = 1.278637452005171
$ws.Cells.Item(10, 14).Value = 1.647356370700479

$ws.Cells.Item(11, 2).Value = 0.4451452708037493
$ws.Cells.Item(11, 3).Value = 0.05006204158269156
$ws.Cells.Item(11, 4).Value = 0.03456937105467262
$ws.Cells.Item(11, 6).Value = 0.9656973931339081
$ws.Cells.Item(11, 7).Value = 0.8093125455051506
$ws.Cells.Item(11, 8).Value = 0.8318359343193151
$ws.Cells.Item(11, 11).Value = 0.4049094492698089
$ws.Cells.Item(11, 13).Value = 1.362696580864935
$ws.Cells.Item(11, 14).Value = 1.637833022364504

$ws.Cells.Item(12, 2).Value = 0.4547780446199852
$ws.Cells.Item(12, 3).Value = 0.05101048768428029
$ws.Cells.Item(12, 4).Value = 0.03484748050702535
$ws.Cells.Item(12, 6).Value = 0.9696560145661124
$ws.Cells.Item(12, 7).Value = 0.8128015748975344
$ws.Cells.Item(12, 8).Value = 0.8324154906566577
$ws.Cells.Item(12, 11).Value = 0.4145585492597661
$ws.Cells.Item(12, 13).Value = 1.3946488791576
$ws.Cells.Item(12, 14).Value = 1.634317038578189

$ws.Cells.Item(13, 2).Value = 0.4527028902537609
$ws.Cells.Item(13, 3).Value = 0.05080634102064607
$ws.Cells.Item(13, 4).Value = 0.03478765312213028
$ws.Cells.Item(13, 6).Value = 0.9687998353013683
$ws.Cells.Item(13, 7).Value = 0.8120468744887432
$ws.Cells.Item(13, 8).Value = 0.8322881940191706
$ws.Cells.Item(13, 11).Value = 0.4124802723470111
$ws.Cells.Item(13, 13).Value = 1.387761882820485
$ws.Cells.Item(13, 14).Value = 1.635070247108807

$ws.Cells.Item(14, 2).Value = 0.4459375102281911
$ws.Cells.Item(14, 3).Value = 0.05014012401730383
$ws.Cells.Item(14, 4).Value = 0.0345922821026079
$ws.Cells.Item(14, 6).Value = 0.9660214407153234
$ws.Cells.Item(14, 7).Value = 0.8095981145621067
$ws.Cells.Item(14, 8).Value = 0.8318824978635604
$ws.Cells.Item(14, 11).Value = 0.405703209220718
$ws.Cells.Item(14, 13).Value = 1.365322849361377
$ws.Cells.Item(14, 14).Value = 1.637541948583262

$ws.Cells.Item(15, 2).Value = 0.4417951768134856
$ws.Cells.Item(15, 3).Value = 0.04973170139452066
$ws.Cells.Item(15, 4).Value = 0.03447241164374759
$ws.Cells.Item(15, 6).Value = 0.9643301875311181
$ws.Cells.Item(15, 7).Value = 0.8081077634380307
$ws.Cells.Item(15, 8).Value = 0.8316412543800737
$ws.Cells.Item(15, 11).Value = 0.4015525670527325
$ws.Cells.Item(15, 13).Value = 1.351594263080585
$ws.Cells.Item(15, 14).Value = 1.639067706880013

$ws.Cells.Item(16, 2).Value = 0.4180827688831812
$ws.Cells.Item(16, 3).Value = 0.04738592581635714
$ws.Cells.Item(16, 4).Value = 0.03378243339876974
$ws.Cells.Item(16, 6).Value = 0.9548004930713176
$ws.Cells.Item(16, 7).Value = 0.7997138104182824
$ws.Cells.Item(16, 8).Value = 0.830369967254498
$ws.Cells.Item(16, 11).Value = 0.3777748781845958
$ws.Cells.Item(16, 13).Value = 1.273160509059991
$ws.Cells.Item(16, 14).Value = 1.647991362278908

$ws.Cells.Item(17, 2).Value = 0.4035590950264805
$ws.Cells.Item(17, 3).Value = 0.04594219199323391
$ws.Cells.Item(17, 4).Value = 0.03335645226948003
$ws.Cells.Item(17, 6).Value = 0.949098868842924
$ws.Cells.Item(17, 7).Value = 0.7946949598590436
$ws.Cells.Item(17, 8).Value = 0.8296892292510734
$ws.Cells.Item(17, 11).Value = 0.3631953899459006
$ws.Cells.Item(17, 13).Value = 1.22525142959384
$ws.Cells.Item(17, 14).Value = 1.6536262211267

$ws.Cells.Item(18, 2).Value = 0.3952140577759451
$ws.Cells.Item(18, 3).Value = 0.04511008641080139
$ws.Cells.Item(18, 4).Value = 0.03311045379273025
$ws.Cells.Item(18, 6).Value = 0.9458724849509963
$ws.Cells.Item(18, 7).Value = 0.7918561193228584
$ws.Cells.Item(18, 8).Value = 0.8293341049032392
$ws.Cells.Item(18, 11).Value = 0.3548124546300357
$ws.Cells.Item(18, 13).Value = 1.197769298738194
$ws.Cells.Item(18, 14).Value = 1.656926115542383

$ws.Cells.Item(19, 2).Value = 0.3923900572998775
$ws.Cells.Item(19, 3).Value = 0.04482805738992113
$ws.Cells.Item(19, 4).Value = 0.03302699440436641
$ws.Cells.Item(19, 6).Value = 0.9447891891763334
$ws.Cells.Item(19, 7).Value = 0.7909031475950314
$ws.Cells.Item(19, 8).Value = 0.8292201187632742
$ws.Cells.Item(19, 11).Value = 0.3519746271132362
$ws.Cells.Item(19, 13).Value = 1.188476794104588
$ws.Cells.Item(19, 14).Value = 1.658053504598193

$ws.Cells.Item(20, 2).Value = 0.4051042775038525
$ws.Cells.Item(20, 3).Value = 0.04609605693977414
$ws.Cells.Item(20, 4).Value = 0.03340190084234251
$ws.Cells.Item(20, 6).Value = 0.9497003253660807
$ws.Cells.Item(20, 7).Value = 0.7952242681968755
$ws.Cells.Item(20, 8).Value = 0.8297579253615908
$ws.Cells.Item(20, 11).Value = 0.3647471137100808
$ws.Cells.Item(20, 13).Value = 1.230343731499929
$ws.Cells.Item(20, 14).Value = 1.653020285330101

$ws.Cells.Item(21, 2).Value = 0.4479243215287738
$ws.Cells.Item(21, 3).Value = 0.0503358801554441
$ws.Cells.Item(21, 4).Value = 0.03464970904276043
$ws.Cells.Item(21, 6).Value = 0.9668353144008393
$ws.Cells.Item(21, 7).Value = 0.810315376559771
$ws.Cells.Item(21, 8).Value = 0.832000148158258
$ws.Cells.Item(21, 11).Value = 0.4076936923102323
$ws.Cells.Item(21, 13).Value = 1.371910396015082
$ws.Cells.Item(21, 14).Value = 1.63681349653811

$ws.Cells.Item(22, 2).Value = 0.4759841345143059
$ws.Cells.Item(22, 3).Value = 0.05309141551904872
$ws.Cells.Item(22, 4).Value = 0.03545628940837275
$ws.Cells.Item(22, 6).Value = 0.9785079639866439
$ws.Cells.Item(22, 7).Value = 0.8206069853408309
$ws.Cells.Item(22, 8).Value = 0.8337903248114031
$ws.Cells.Item(22, 11).Value = 0.4357846166890056
$ws.Cells.Item(22, 13).Value = 1.465140809049089
$ws.Cells.Item(22, 14).Value = 1.626747926391204

$ws.Cells.Item(23, 2).Value = 0.4610013749589257
$ws.Cells.Item(23, 3).Value = 0.05162215755679256
$ws.Cells.Item(23, 4).Value = 0.03502662751023422
$ws.Cells.Item(23, 6).Value = 0.9722346055842621
$ws.Cells.Item(23, 7).Value = 0.8150748160542349
$ws.Cells.Item(23, 8).Value = 0.8328051362899629
$ws.Cells.Item(23, 11).Value = 0.4207899717383725
$ws.Cells.Item(23, 13).Value = 1.415314690538509
$ws.Cells.Item(23, 14).Value = 1.632071825415458

$ws.Cells.Item(24, 2).Value = 0.4044056854805831
$ws.Cells.Item(24, 3).Value = 0.04602650108606099
$ws.Cells.Item(24, 4).Value = 0.03338135695498323
$ws.Cells.Item(24, 6).Value = 0.9494282463049899
$ws.Cells.Item(24, 7).Value = 0.7949848228889635
$ws.Cells.Item(24, 8).Value = 0.8297267549720573
$ws.Cells.Item(24, 11).Value = 0.3640455824672983
$ws.Cells.Item(24, 13).Value = 1.228041310997895
$ws.Cells.Item(24, 14).Value = 1.653294040985863

$ws.Cells.Item(25, 2).Value = 0.3436877254603985
$ws.Cells.Item(25, 3).Value = 0.03992548624762549
$ws.Cells.Item(25, 4).Value = 0.03156908162210925
$ws.Cells.Item(25, 6).Value = 0.9268546352191436
$ws.Cells.Item(25, 7).Value = 0.7751431425863871
$ws.Cells.Item(25, 8).Value = 0.8277982722988213
$ws.Cells.Item(25, 11).Value = 0.3029456556859884
$ws.Cells.Item(25, 13).Value = 1.028833229057994
$ws.Cells.Item(25, 14).Value = 1.678435218975068
